$wb = $excel.ActiveWorkbook

# --- Sheet1: update the "Conversión del día" text cell (A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("✅ 1000 Bs = 3.7 = 14225.93 pesos", "✅ 1000 Bs = 3.67 = 14020.54 pesos")
$text = $text.Replace("✅ 14225.93 pesos = 3.7 = 931.15 Bs", "✅ 14020.54 pesos = 3.66 = 957.31 Bs")
$cell.Value = $text

# --- Sheet2 ("tasas"): update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 272.6
$ws2.Range("O10").Value = 3822
$ws2.Range("N12").Value = 3830
$ws2.Range("O12").Value = 261.51
